$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 258 & 259: add column G (YouTube Channel link) ---
$ws.Hyperlinks.Add($ws.Range("G258"), "https://www.youtube.com/@805Webcams/streams", [Type]::Missing, [Type]::Missing, "https://www.youtube.com/@805Webcams/streams")
$ws.Range("G258").Value = "Morro Bay Jetty Webcam - 805 Webcams in California"

$ws.Hyperlinks.Add($ws.Range("G259"), "https://www.youtube.com/@805Webcams/streams", [Type]::Missing, [Type]::Missing, "https://www.youtube.com/@805Webcams/streams")
$ws.Range("G259").Value = "Morro Bay Jetty Webcam - 805 Webcams in California"

# --- New row 260 ---
$ws.Range("A260").Value = "LIVE, CHIMNEY, ROCK"
$ws.Range("B260").Value = "35.36988113515808, -120.86660169717162"
$ws.Range("C260").Value = "Morro Bay Jetty Webcam - 805 Webcams in California"
$ws.Range("D260").Value = "CA"
$ws.Range("E260").Value = "USA"
$ws.Range("F260").Value = "g1F2ktr4e10"
$ws.Hyperlinks.Add($ws.Range("G260"), "https://www.youtube.com/@805Webcams/streams", [Type]::Missing, [Type]::Missing, "https://www.youtube.com/@805Webcams/streams")
$ws.Range("G260").Value = "Morro Bay Jetty Webcam - 805 Webcams in California"

# --- New row 261 ---
$ws.Range("A261").Value = "LIVE, STREET, PARK"
$ws.Range("B261").Value = "42.352590061743065, -71.06684578819161"
$ws.Range("C261").Value = "Boston Common / Boylston St. Live Cam"
$ws.Range("D261").Value = "MA"
$ws.Range("E261").Value = "USA"
$ws.Range("F261").Value = "4nYY5p6ClUU"
$ws.Hyperlinks.Add($ws.Range("G261"), "https://www.youtube.com/@BostonAndMaineLive/streams", [Type]::Missing, [Type]::Missing, "https://www.youtube.com/@BostonAndMaineLive/streams")
$ws.Range("G261").Value = "(170) Boston and Maine Live - YouTube"

# --- New row 262 ---
$ws.Range("A262").Value = "LIVE, SEA, BEACH"
$ws.Range("B262").Value = "43.13209756110367, -70.63838259811487"
$ws.Range("C262").Value = "York Harbor Beach"
$ws.Range("D262").Value = "ME"
$ws.Range("E262").Value = "USA"
$ws.Range("F262").Value = "catvjIWNrZg"
$ws.Hyperlinks.Add($ws.Range("G262"), "https://www.youtube.com/@BostonAndMaineLive/streams", [Type]::Missing, [Type]::Missing, "https://www.youtube.com/@BostonAndMaineLive/streams")
$ws.Range("G262").Value = "(170) Boston and Maine Live - YouTube"

# --- New row 263 ---
$ws.Range("A263").Value = "LIVE, SEA, BEACH"
$ws.Range("B263").Value = "43.159514229404174, -70.62047217223258"
$ws.Range("C263").Value = "York Beach, Maine US - Anchorage Inn"
$ws.Range("D263").Value = "ME"
$ws.Range("E263").Value = "USA"
$ws.Range("F263").Value = "bnUgt0gl-ds"
$ws.Hyperlinks.Add($ws.Range("G263"), "https://www.youtube.com/@BostonAndMaineLive/streams", [Type]::Missing, [Type]::Missing, "https://www.youtube.com/@BostonAndMaineLive/streams")
$ws.Range("G263").Value = "(170) Boston and Maine Live - YouTube"

# --- New row 264 ---
$ws.Range("A264").Value = "LIVE, RAIL, TRAIN, BRIDGE"
$ws.Range("B264").Value = "42.12949635554303, -72.74572011287621"
$ws.Range("C264").Value = "Westfield, Massachusetts USA - LIVE Railfan Cam"
$ws.Range("D264").Value = "MA"
$ws.Range("E264").Value = "USA"
$ws.Range("F264").Value = "njLfO0eQ_Ug"
$ws.Hyperlinks.Add($ws.Range("G264"), "https://www.youtube.com/@BostonAndMaineLive/streams", [Type]::Missing, [Type]::Missing, "https://www.youtube.com/@BostonAndMaineLive/streams")
$ws.Range("G264").Value = "(170) Boston and Maine Live - YouTube"

# Apply the same border/style used by the rest of column A/E for these new data rows
$ws.Range("A260:A264").Style = $ws.Range("A258").Style
$ws.Range("E260:E264").Style = $ws.Range("E258").Style

# Move the active selection to the next empty row, like Excel leaves after pasting new rows
$ws.Range("A265").Select()
